$wb = $excel.ActiveWorkbook

# Duplicate the "Slovakia" sheet (closest template) to the end of the tab
# strip, then rename/re-point it to become the new "Italy" sheet.
$slovakia = $wb.Worksheets.Item("Slovakia")
$slovakia.Copy($null, $slovakia)
$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"

# Before changing contents, select the whole sheet on Slovakia so it is
# left with a "select all" selection state (matches a tab switch away).
$slovakia.Activate()
$slovakia.Cells.Select()

# Populate the new sheet's market name / part-number cells.
$italy.Range("B2").Value = "Italy Market"
$italy.Range("B4").ClearFormats()
$italy.Range("B4").Value = "NGC-3145/T2160"

# Make Italy the active tab, with the cursor parked on E21.
$italy.Activate()
$italy.Range("E21").Select()
